$wb = $excel.ActiveWorkbook

# hunk0: ALC row9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 9058.416999999999
$ws.Range("I9").Value = 16869.834
$ws.Range("K9").Value = 16869.834
$ws.Range("M9").Value = -16700.834

# hunk1: ALC row34
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2750
$ws.Range("I34").Value = 2750
$ws.Range("K34").Value = 2750
$ws.Range("M34").Value = -2547

# hunk2: ALC row36
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 2750
$ws.Range("I36").Value = 2750
$ws.Range("K36").Value = 2750
$ws.Range("M36").Value = -2035

# hunk3: ALC row39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 395.66666
$ws.Range("J39").Value = 1025
$ws.Range("L39").Value = 3075
$ws.Range("N39").Value = -3667

# hunk4: ALC row131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 7500

# hunk5: ARM row97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2641.8235
$ws.Range("I97").Value = 2641.8235
$ws.Range("K97").Value = 2641.8235
$ws.Range("M97").Value = -2145.8235

# hunk6: ARM row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3021.9678
$ws.Range("I122").Value = 2065.1
$ws.Range("K122").Value = 6195.299999999999
$ws.Range("M122").Value = -3745.299999999999

# hunk7: ARM row124
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 26815.2
$ws.Range("J124").Value = 26815.2
$ws.Range("L124").Value = 26815.2
$ws.Range("N124").Value = -36635.2

# hunk8: ARM row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 28573832
$ws.Range("I132").Value = 1884.6897
$ws.Range("J132").Value = 166671580
$ws.Range("K132").Value = 5654.0691
$ws.Range("L132").Value = 500014740
$ws.Range("M132").Value = -3124.0691
$ws.Range("N132").Value = -500019800

# hunk9: BSM row26
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 31838.666
$ws.Range("I26").Value = 17498.5
$ws.Range("K26").Value = 17498.5
$ws.Range("M26").Value = -17206.5

# hunk10: BSM row106
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 138618.2
$ws.Range("J106").Value = 138618.2
$ws.Range("L106").Value = 138618.2
$ws.Range("N106").Value = -141142.2

# hunk11: BSM row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5407.3335
$ws.Range("J107").Value = 4694.3335
$ws.Range("L107").Value = 4694.3335
$ws.Range("N107").Value = -8534.333500000001

# hunk12: CRP row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30308170
$ws.Range("J31").Value = 125010300
$ws.Range("L31").Value = 125010300
$ws.Range("N31").Value = -125010890

# hunk13: CRP row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 30308170
$ws.Range("J34").Value = 125010300
$ws.Range("L34").Value = 125010300
$ws.Range("N34").Value = -125010704

# hunk14: CRP row68
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 39221.25
$ws.Range("J68").Value = 39221.25
$ws.Range("L68").Value = 39221.25
$ws.Range("N68").Value = -40719.25

# hunk15: CRP row71
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 39221.25
$ws.Range("J71").Value = 39221.25
$ws.Range("L71").Value = 117663.75
$ws.Range("N71").Value = -125151.75

# hunk16: CRP row74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 40157
$ws.Range("J74").Value = 40157
$ws.Range("L74").Value = 40157
$ws.Range("N74").Value = -41905

# hunk17: CRP row77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 40157
$ws.Range("J77").Value = 40157
$ws.Range("L77").Value = 120471
$ws.Range("N77").Value = -129207

# hunk18: CRP row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8256.723
$ws.Range("I99").Value = 6777.2
$ws.Range("K99").Value = 6777.2
$ws.Range("M99").Value = -5279.2

# hunk19: CRP row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2303.923
$ws.Range("I122").Value = 2111.5557
$ws.Range("J122").Value = 2736.75
$ws.Range("K122").Value = 6334.6671
$ws.Range("L122").Value = 8210.25
$ws.Range("M122").Value = -3884.6671
$ws.Range("N122").Value = -13110.25

# hunk20: CRP row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 8256.723
$ws.Range("I126").Value = 6777.2
$ws.Range("K126").Value = 20331.6
$ws.Range("M126").Value = -17861.6

# hunk21: CRP row141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 285328.38
$ws.Range("J141").Value = 364971.9
$ws.Range("L141").Value = 364971.9
$ws.Range("N141").Value = -375331.9

# hunk22: CUL row60
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 990.3077
$ws.Range("I60").Value = 252.27272
$ws.Range("K60").Value = 756.81816
$ws.Range("M60").Value = -505.81816

# hunk23: CUL row116
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 2500
$ws.Range("I116").Value = 500
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 1500
$ws.Range("L116").Value = 9000
$ws.Range("M116").Value = 1942
$ws.Range("N116").Value = -15884

# hunk24: CUL row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 38374.184
$ws.Range("I131").Value = 149481.42
$ws.Range("K131").Value = 448444.26
$ws.Range("M131").Value = -443404.26

# hunk25: GSM row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1884.6487
$ws.Range("I102").Value = 1350.36
$ws.Range("J102").Value = 2997.75
$ws.Range("K102").Value = 1350.36
$ws.Range("L102").Value = 2997.75
$ws.Range("M102").Value = 271.6400000000001
$ws.Range("N102").Value = -6241.75

# hunk26: GSM row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1914.9333
$ws.Range("J122").Value = 1932.6
$ws.Range("L122").Value = 5797.799999999999
$ws.Range("N122").Value = -10697.8

# hunk27: GSM row140 (M140 removed, N140 added)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 90000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 90000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 90000
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -100360

# hunk28: LTW row2 (N2 removed)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 105
$ws.Range("I2").Value = 105
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 105
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 7
$ws.Range("N2").ClearContents()

# hunk29: LTW row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2927.4285
$ws.Range("I7").Value = 2298.8
$ws.Range("K7").Value = 2298.8
$ws.Range("M7").Value = -2186.8

# hunk30: LTW row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1723.1708
$ws.Range("I46").Value = 619.35486
$ws.Range("J46").Value = 5145
$ws.Range("K46").Value = 619.35486
$ws.Range("L46").Value = 5145
$ws.Range("M46").Value = -431.35486
$ws.Range("N46").Value = -5521

# hunk31: LTW row55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 290.76
$ws.Range("J55").Value = 314.63635
$ws.Range("L55").Value = 314.63635
$ws.Range("N55").Value = -660.63635

# hunk32: LTW row100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4299.5
$ws.Range("I100").Value = 2799
$ws.Range("J100").Value = 4799.6665
$ws.Range("K100").Value = 2799
$ws.Range("L100").Value = 4799.6665
$ws.Range("M100").Value = -2258
$ws.Range("N100").Value = -5881.6665

# hunk33: LTW row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2927.4285
$ws.Range("I126").Value = 2298.8
$ws.Range("K126").Value = 6896.400000000001
$ws.Range("M126").Value = -4426.400000000001

# hunk34: WVR row52 (N52 removed, M52 unchanged)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 10660.5
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# hunk35: WVR row70
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 55666.332
$ws.Range("I70").Value = 40999.5
$ws.Range("K70").Value = 40999.5
$ws.Range("M70").Value = -40684.5

# hunk36: WVR row73
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 55666.332
$ws.Range("I73").Value = 40999.5
$ws.Range("K73").Value = 40999.5
$ws.Range("M73").Value = -39907.5

# hunk37: WVR row96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7618.2
$ws.Range("I96").Value = 6147
$ws.Range("J96").Value = 8248.714
$ws.Range("K96").Value = 6147
$ws.Range("L96").Value = 8248.714
$ws.Range("M96").Value = -4774
$ws.Range("N96").Value = -10994.714

# hunk38: WVR row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 50051988
$ws.Range("I122").Value = 62564010
$ws.Range("J122").Value = 3912.5
$ws.Range("K122").Value = 187692030
$ws.Range("L122").Value = 11737.5
$ws.Range("M122").Value = -187689580
$ws.Range("N122").Value = -16637.5

# hunk39: WVR row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4091.48
$ws.Range("I126").Value = 4212
$ws.Range("K126").Value = 12636
$ws.Range("M126").Value = -10166

# hunk40: WVR row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2288.318
$ws.Range("I136").Value = 1632.3334
$ws.Range("K136").Value = 4897.0002
$ws.Range("M136").Value = -2347.0002
